# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Rules sheet, cell B11 ("R40") is retyped to the text "1" (still a text
# cell, not a number) - the 4th rule name in the decision table becomes "1".
#
# A plain $ws.Range("B11").Value = "1" would let Excel auto-detect the
# literal as a number, which changes the cell's type/format. To keep it a
# genuine text value (matching the original "R40" text cell), build it as a
# text formula and paste back as a value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial($xlPasteValues)
